$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.517.24'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('D3').Value = '1.677.03'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''219.63'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').Value = '''0.5318'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('D7').Value = '''1.002'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '''0.2692'
$ws.Range('E8').Value = '  +3.11%  '
$ws.Range('D9').Value = '''0.06403'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').Value = '''21.79'
$ws.Range('E10').Value = '  +4.57%  '
$ws.Range('D11').Value = '''0.07805'
$ws.Range('D12').Value = '1.738.17'
$ws.Range('E12').Value = '  +5.55%  '
$ws.Range('D13').Value = '''4.504'
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').Value = '''0.5582'
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = '0.0₅8346'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '''65.75'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = '26.534.68'
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '''4.786'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '''192.63'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '''10.33'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = '''6.321'
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '''0.1275'
$ws.Range('E24').Value = '  +5.08%  '
$ws.Range('D25').Value = '''139.16'
$ws.Range('E25').Value = '  -5.03%  '
$ws.Range('D26').Value = '''7.425'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').Value = '''16.28'
$ws.Range('E27').Value = '  +2.47%  '
$ws.Range('D28').Value = '''1.433'
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('D29').Value = '''0.06318'
$ws.Range('E29').Value = '  +6.98%  '
$ws.Range('D30').Value = '''1.286'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = '''3.605'
$ws.Range('E31').Value = '  +5.13%  '
$ws.Range('D32').Value = '''3.441'
$ws.Range('D33').Value = '''1.692'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('D34').Value = '''1.014'
$ws.Range('E34').Value = '  +2.62%  '
$ws.Range('D35').Value = '''0.6145'
$ws.Range('E35').Value = '  +8.34%  '
$ws.Range('D36').Value = '''2.424'
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('D37').Value = '''2.786'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').Value = '''0.01631'
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('D39').Value = '''6.092'
$ws.Range('E39').Value = '  +4.69%  '
$ws.Range('D40').Value = '1.094.20'
$ws.Range('E40').Value = '  +6.08%  '
$ws.Range('D41').Value = '''0.8620'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = '''100.61'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').Value = '1.824.07'
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₈113'
$ws.Range('E45').Value = '  +3.96%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''58.60'
$ws.Range('E46').Value = '  +4.45%  '
$ws.Range('D47').Value = '''8.192'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('D48').Value = '''0.9988'
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('D49').Value = '''1.503'
$ws.Range('E49').Value = '  +8.54%  '
$ws.Range('D50').Value = '''0.05198'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('E51').Value = '  +1.45%  '
